$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 500
$ws.Range("B4").Value = 300
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = "[ 2.81672028e-02 -6.64177805e-02  1.21331494e-01  3.18853641e-01`n  5.14147939e-01  1.84955188e-01  9.82713578e-01 -2.85932855e-02`n  8.53394519e-01  3.25549233e-01 -9.41669520e-02  2.18408589e-01`n  7.43831425e-02 -1.93778443e-01 -2.83944315e-01  2.26082456e-04`n  3.19551617e-01 -1.51715657e-01]"
$ws.Range("F4").Value = 0.9576059850374065
$ws.Range("G4").Value = 0.9306930693069307
$ws.Range("H4").Value = 0.9251101321585903
$ws.Range("I4").Value = 0.8771929824561403

$ws.Rows.Item(4).AutoFit()
